$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label for the new "Total" column
$ws.Range("Q3").Value = "Total"
$ws.Range("Q3").Style = $ws.Range("C3").Style

# Row totals: Q4:Q16 = SUM(E:P) for that row, using the same number style as the data cells
for ($r = 4; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 17)
    $cell.Formula = "=SUM(E" + $r + ":P" + $r + ")"
    $cell.Style = $ws.Cells.Item($r, 16).Style
}

# Update the active selection to match the recorded view state
$ws.Range("U5").Select()
